# recovery-factor.xlsx update
# - Extends the data table from 10 observations (rows 2-11) to 18 observations
#   (rows 2-19), recording flat/placeholder outputs (B=1000, C=0.1, D=100, E=100)
#   instead of the old varying numbers.
# - Adds two new analysis columns: G = ABS(D-E), H = running smoothing of G.
# - Re-points the totals row (now row 20) and the chart source ranges to match.
# - K2 (the smoothing constant) goes from 0.3 to 1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Make room: push the old "totals" row (row 12) down to row 20 by
#    inserting 8 fresh rows above it.
# ---------------------------------------------------------------------------
$ws.Range("A12:K19").EntireRow.Insert()

# ---------------------------------------------------------------------------
# 2. Smoothing constant K2: 0.3 -> 1
# ---------------------------------------------------------------------------
$ws.Range("K2").Value = 1

# ---------------------------------------------------------------------------
# 3. Restore column B ("N") for the original rows to a flat 1000 (rows
#    6, 8, 10 and 11 had been 1100 / 1220 / 750 / 880) - row 4 (1020) is
#    left untouched.
# ---------------------------------------------------------------------------
$ws.Range("B6").Value = 1000
$ws.Range("B8").Value = 1000
$ws.Range("B10").Value = 1000
$ws.Range("B11").Value = 1000

# ---------------------------------------------------------------------------
# 4. Column D ("Y"): rows 4,6,8,9,10,11 now hold the literal value 100
#    (no longer driven by the B*C formula) - same "record outputs" flattening
#    the rest of the table gets.
# ---------------------------------------------------------------------------
$ws.Range("D4").Value = 100
$ws.Range("D6").Value = 100
$ws.Range("D8").Value = 100
$ws.Range("D9").Value = 100
$ws.Range("D10").Value = 100
$ws.Range("D11").Value = 100

# ---------------------------------------------------------------------------
# 5. New column G ("R^2" helper) = ABS(D-E) for every data row, rows 2-19.
# ---------------------------------------------------------------------------
$ws.Range("G2:G19").Formula = '=ABS(D2-E2)'

# ---------------------------------------------------------------------------
# 6. New column H: a running pairwise average of G, starting at row 3.
#    H3 is special (root-mean-square style average of the first two G's),
#    every other row from H4 onward is a plain average of the current and
#    previous G.
# ---------------------------------------------------------------------------
$ws.Range("H3").Formula = '=(G2^2+G3^2)/(G2+G3)'
$ws.Range("H4:H19").Formula = '=(G4+G3)/2'

# ---------------------------------------------------------------------------
# 7. Populate the 8 newly-inserted rows (12-19 => observations 11-18) with
#    the same flat placeholder values/formulas as the rest of the table.
# ---------------------------------------------------------------------------
for ($i = 12; $i -le 19; $i++) {
    $prev = $i - 1
    $ws.Cells.Item($i, 1).Value = $i - 1
    $ws.Cells.Item($i, 2).Value = 1000
    $ws.Cells.Item($i, 3).Value = 0.1
    $ws.Cells.Item($i, 4).Value = 100

    $eFormula = '=ROUNDUP((SUM(D$2:D' + $prev + ') - SUM($E$2:E' + $prev + '))*$K$2 + D' + $i + ', 0)'
    $ws.Cells.Item($i, 5).Formula = $eFormula
}

# F column ((D-E)^2) for the new rows - rows 2-11 already have it.
$ws.Range("F12:F19").Formula = '=(D12-E12)^2'

# ---------------------------------------------------------------------------
# 8. Fix up the totals row, now row 20: it only ever summed the first 12
#    observations (D2:D13 / E2:E13 / F2:F13), plus a new H20 = SUM(H3:H13).
# ---------------------------------------------------------------------------
$ws.Range("D20:E20").Formula = '=SUM(D2:D13)'
$ws.Range("F20").Formula = '=SUM(F2:F13)'
$ws.Range("H20").Formula = '=SUM(H3:H13)'

# ---------------------------------------------------------------------------
# 9. Update the line chart so both series (D = "STAT"/Y and E = OBS) plot
#    the full 18-row range instead of the original 10.
# ---------------------------------------------------------------------------
$co = $ws.ChartObjects().Item(1)
$chart = $co.Chart
$s1 = $chart.SeriesCollection().Item(1)
$s1.Formula = '=SERIES(Sheet1!$D$1,Sheet1!$A$2:$A$19,Sheet1!$D$2:$D$19,1)'
$s2 = $chart.SeriesCollection().Item(2)
$s2.Formula = '=SERIES(Sheet1!$E$1,Sheet1!$A$2:$A$19,Sheet1!$E$2:$E$19,2)'

# ---------------------------------------------------------------------------
# 10. Minor cosmetic touches matching the authored workbook.
# ---------------------------------------------------------------------------
$ws.Range("D2").Select()
